$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 587.5714
$ws.Range("I80").Value = 699.375
$ws.Range("J80").Value = 438.5
$ws.Range("K80").Value = 2098.125
$ws.Range("L80").Value = 1315.5
$ws.Range("M80").Value = -1100.125
$ws.Range("N80").Value = -3311.5
$ws.Range("H83").Value = 587.5714
$ws.Range("I83").Value = 699.375
$ws.Range("J83").Value = 438.5
$ws.Range("K83").Value = 6294.375
$ws.Range("L83").Value = 3946.5
$ws.Range("M83").Value = -1302.375
$ws.Range("N83").Value = -13930.5
$ws.Range("H86").Value = 3100457.5
$ws.Range("I86").Value = 4832
$ws.Range("J86").Value = 4788980.5
$ws.Range("K86").Value = 4832
$ws.Range("L86").Value = 4788980.5
$ws.Range("M86").Value = -3709
$ws.Range("N86").Value = -4791226.5
$ws.Range("H89").Value = 3100457.5
$ws.Range("I89").Value = 4832
$ws.Range("J89").Value = 4788980.5
$ws.Range("K89").Value = 24160
$ws.Range("L89").Value = 23944902.5
$ws.Range("M89").Value = -18544
$ws.Range("N89").Value = -23956134.5
$ws.Range("H98").Value = 1835.561
$ws.Range("I98").Value = 1190.7354
$ws.Range("K98").Value = 1190.7354
$ws.Range("M98").Value = 307.2646
$ws.Range("H101").Value = 1014.5
$ws.Range("I101").Value = 1016.5714
$ws.Range("K101").Value = 3049.7142
$ws.Range("M101").Value = -1427.7142
$ws.Range("H107").Value = 37382.926
$ws.Range("I107").Value = 38813.04
$ws.Range("K107").Value = 38813.04
$ws.Range("M107").Value = -36893.04
$ws.Range("H122").Value = 1835.561
$ws.Range("I122").Value = 1190.7354
$ws.Range("K122").Value = 3572.2062
$ws.Range("M122").Value = -1122.2062
$ws.Range("H125").Value = 7939462
$ws.Range("I125").Value = 2185.8
$ws.Range("K125").Value = 19672.2
$ws.Range("M125").Value = -17212.2
$ws.Range("H129").Value = 2277.4092
$ws.Range("I129").Value = 1074.8572
$ws.Range("J129").Value = 2838.6
$ws.Range("K129").Value = 3224.5716
$ws.Range("L129").Value = 8515.799999999999
$ws.Range("M129").Value = 1775.4284
$ws.Range("N129").Value = -18515.8
$ws.Range("H137").Value = 647783.6
$ws.Range("I137").Value = 557616.75
$ws.Range("J137").Value = 772630
$ws.Range("K137").Value = 1672850.25
$ws.Range("L137").Value = 2317890
$ws.Range("M137").Value = -1670300.25
$ws.Range("N137").Value = -2322990
$ws.Range("H141").Value = 3112.9092
$ws.Range("I141").Value = 3174.2
$ws.Range("K141").Value = 9522.599999999999
$ws.Range("M141").Value = -4342.599999999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3994.5334
$ws.Range("I61").Value = 2175.5557
$ws.Range("K61").Value = 2175.5557
$ws.Range("M61").Value = -1963.5557
$ws.Range("H74").Value = 1417.5
$ws.Range("I74").Value = 1302.9412
$ws.Range("J74").Value = 2066.6667
$ws.Range("K74").Value = 1302.9412
$ws.Range("L74").Value = 2066.6667
$ws.Range("M74").Value = -428.9412
$ws.Range("N74").Value = -3814.6667
$ws.Range("H77").Value = 1417.5
$ws.Range("I77").Value = 1302.9412
$ws.Range("J77").Value = 2066.6667
$ws.Range("K77").Value = 6514.706
$ws.Range("L77").Value = 10333.3335
$ws.Range("M77").Value = -2146.706
$ws.Range("N77").Value = -19069.3335
$ws.Range("H102").Value = 2233
$ws.Range("I102").Value = 2233
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2233
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -611
$ws.Range("H136").Value = 3994.5334
$ws.Range("I136").Value = 2175.5557
$ws.Range("K136").Value = 6526.6671
$ws.Range("M136").Value = -3976.6671

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 236
$ws.Range("I22").Value = 183.28572
$ws.Range("K22").Value = 183.28572
$ws.Range("M22").Value = -10.28572
$ws.Range("H86").Value = 1701835.8
$ws.Range("I86").Value = 2127069.5
$ws.Range("K86").Value = 2127069.5
$ws.Range("M86").Value = -2125946.5
$ws.Range("H89").Value = 1701835.8
$ws.Range("I89").Value = 2127069.5
$ws.Range("K89").Value = 10635347.5
$ws.Range("M89").Value = -10629731.5
$ws.Range("H107").Value = 345630.2
$ws.Range("I107").Value = 897.1739
$ws.Range("K107").Value = 897.1739
$ws.Range("M107").Value = 1022.8261
$ws.Range("H134").Value = 33702.395
$ws.Range("I134").Value = 2333.1365
$ws.Range("K134").Value = 6999.4095
$ws.Range("M134").Value = -4464.4095

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25329.273
$ws.Range("I31").Value = 1544.1333
$ws.Range("K31").Value = 1544.1333
$ws.Range("M31").Value = -1249.1333
$ws.Range("H34").Value = 25329.273
$ws.Range("I34").Value = 1544.1333
$ws.Range("K34").Value = 1544.1333
$ws.Range("M34").Value = -1342.1333
$ws.Range("H58").Value = 6432.6
$ws.Range("I58").Value = 5467.091
$ws.Range("J58").Value = 7612.6665
$ws.Range("K58").Value = 5467.091
$ws.Range("L58").Value = 7612.6665
$ws.Range("M58").Value = -5264.091
$ws.Range("N58").Value = -8018.6665
$ws.Range("H62").Value = 5199.8
$ws.Range("I62").Value = 1750
$ws.Range("J62").Value = 7499.6665
$ws.Range("K62").Value = 1750
$ws.Range("L62").Value = 7499.6665
$ws.Range("M62").Value = -1126
$ws.Range("N62").Value = -8747.666499999999
$ws.Range("H65").Value = 5199.8
$ws.Range("I65").Value = 1750
$ws.Range("J65").Value = 7499.6665
$ws.Range("K65").Value = 8750
$ws.Range("L65").Value = 37498.3325
$ws.Range("M65").Value = -5630
$ws.Range("N65").Value = -43738.3325
$ws.Range("H122").Value = 4683.3335
$ws.Range("I122").Value = 2831.375
$ws.Range("J122").Value = 6799.857
$ws.Range("K122").Value = 8494.125
$ws.Range("L122").Value = 20399.571
$ws.Range("M122").Value = -6044.125
$ws.Range("N122").Value = -25299.571
$ws.Range("H132").Value = 4347.9697
$ws.Range("I132").Value = 4279.8096
$ws.Range("J132").Value = 4467.25
$ws.Range("K132").Value = 12839.4288
$ws.Range("L132").Value = 13401.75
$ws.Range("M132").Value = -10309.4288
$ws.Range("N132").Value = -18461.75
$ws.Range("H134").Value = 912275.4399999999
$ws.Range("I134").Value = 558502.6
$ws.Range("K134").Value = 1675507.8
$ws.Range("M134").Value = -1672972.8
$ws.Range("H136").Value = 6432.6
$ws.Range("I136").Value = 5467.091
$ws.Range("J136").Value = 7612.6665
$ws.Range("K136").Value = 16401.273
$ws.Range("L136").Value = 22837.9995
$ws.Range("M136").Value = -13851.273
$ws.Range("N136").Value = -27937.9995

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1818316
$ws.Range("I11").Value = 1818316
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 5454948
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -5454808
$ws.Range("H12").Value = 146.78572
$ws.Range("J12").Value = 158
$ws.Range("L12").Value = 474
$ws.Range("N12").Value = -820
$ws.Range("H26").Value = 604.2857
$ws.Range("I26").Value = 446
$ws.Range("J26").Value = 1000
$ws.Range("K26").Value = 1338
$ws.Range("L26").Value = 3000
$ws.Range("M26").Value = -1050
$ws.Range("N26").Value = -3576
$ws.Range("H32").Value = 9000700
$ws.Range("I32").Value = 30000000
$ws.Range("J32").Value = 2000933.4
$ws.Range("K32").Value = 90000000
$ws.Range("L32").Value = 6002800.199999999
$ws.Range("M32").Value = -89999717
$ws.Range("N32").Value = -6003366.199999999
$ws.Range("H46").Value = 1725
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1725
$ws.Range("K46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("M46").Value = 5175
$ws.Range("N46").Value = -5357
$ws.Range("H86").Value = 362.25
$ws.Range("I86").Value = 275
$ws.Range("J86").Value = 449.5
$ws.Range("K86").Value = 825
$ws.Range("L86").Value = 1348.5
$ws.Range("M86").Value = 361
$ws.Range("N86").Value = -3720.5
$ws.Range("H89").Value = 362.25
$ws.Range("I89").Value = 275
$ws.Range("J89").Value = 449.5
$ws.Range("K89").Value = 2475
$ws.Range("L89").Value = 4045.5
$ws.Range("M89").Value = 3453
$ws.Range("N89").Value = -15901.5
$ws.Range("H131").Value = 15972362
$ws.Range("I131").Value = 47762596
$ws.Range("K131").Value = 143287788
$ws.Range("M131").Value = -143282748

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 776.4
$ws.Range("J107").Value = 1047.3636
$ws.Range("L107").Value = 1047.3636
$ws.Range("N107").Value = -4887.3636
$ws.Range("H122").Value = 3158.9524
$ws.Range("I122").Value = 1839.5
$ws.Range("K122").Value = 5518.5
$ws.Range("M122").Value = -3068.5
$ws.Range("H126").Value = 3107.4614
$ws.Range("I126").Value = 2328.4285
$ws.Range("K126").Value = 6985.2855
$ws.Range("M126").Value = -4515.2855

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 389417.94
$ws.Range("J7").Value = 914137.9399999999
$ws.Range("L7").Value = 914137.9399999999
$ws.Range("N7").Value = -914361.9399999999
$ws.Range("H55").Value = 835.8
$ws.Range("I55").Value = 206.3
$ws.Range("K55").Value = 206.3
$ws.Range("M55").Value = -33.30000000000001
$ws.Range("H126").Value = 389417.94
$ws.Range("J126").Value = 914137.9399999999
$ws.Range("L126").Value = 2742413.82
$ws.Range("N126").Value = -2747353.82

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 7500
$ws.Range("J44").Value = 7500
$ws.Range("L44").Value = 7500
$ws.Range("N44").Value = -8608
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("N99").Value = 0
$ws.Range("H107").Value = 395.56
$ws.Range("I107").Value = 398.83334
$ws.Range("K107").Value = 1196.50002
$ws.Range("M107").Value = 723.4999800000001
$ws.Range("H132").Value = 38190.4
$ws.Range("I132").Value = 4248.5713
$ws.Range("K132").Value = 12745.7139
$ws.Range("M132").Value = -10215.7139
